# Auto-generated Excel COM-interop script: add the 4/17/2020 (serial 43937) data column
# to each of the four sheets in the DC COVID-19 tracking workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet "Overal Stats": new date column AQ (4/17/2020) ---
$ws1 = $wb.Worksheets.Item("Overal Stats")
$ws1.Range("AQ1").Value2 = 43937
$ws1.Range("AQ1").NumberFormat = $ws1.Range("AP1").NumberFormat
$ws1.Range("AQ3").Value2 = 12643
$ws1.Range("AQ4").Value2 = 2476
$ws1.Range("AQ5").Value2 = 86
$ws1.Range("AQ6").Value2 = 573
$ws1.Range("AQ8").Value2 = 98
$ws1.Range("AQ9").Value2 = 443
$ws1.Range("AQ10").Value2 = 204
$ws1.Range("AQ11").Value2 = 239
$ws1.Range("AQ15").Value2 = 74
$ws1.Range("AQ16").Value2 = 41
$ws1.Range("AQ17").Value2 = 33
$ws1.Range("AQ18").Value2 = 119
$ws1.Range("AQ19").Value2 = 119
$ws1.Range("AQ20").Value2 = 502
$ws1.Range("AQ23").Value2 = 74
$ws1.Range("AQ24").Value2 = 57
$ws1.Range("AQ25").Value2 = 17
$ws1.Range("AQ26").Value2 = 139
$ws1.Range("AQ27").Value2 = 196
$ws1.Range("AQ28").Value2 = 559
$ws1.Range("AQ31").Value2 = 20
$ws1.Range("AQ32").Value2 = 20
$ws1.Range("AQ33").Value2 = 0
$ws1.Range("AQ34").Value2 = 149
$ws1.Range("AQ35").Value2 = 168
$ws1.Range("AQ36").Value2 = 53
$ws1.Range("AQ39").Value2 = 65
$ws1.Range("AQ40").Value2 = 22
$ws1.Range("AQ41").Value2 = 43
$ws1.Range("AQ42").Value2 = 512
$ws1.Range("AQ43").Value2 = 534
$ws1.Range("AQ44").Value2 = 43
$ws1.Range("AQ45").Value2 = 1
$ws1.Range("AQ48").Value2 = 13
$ws1.Range("AQ49").Value2 = 12
$ws1.Range("AQ50").Value2 = 0
$ws1.Range("AQ51").Value2 = 63
$ws1.Range("AQ52").Value2 = 75
$ws1.Range("AQ53").Value2 = 46
$ws1.Range("AQ54").Value2 = 1
$ws1.Range("AQ56").Value2 = 7
$ws1.Range("AQ57").Value2 = 6
$ws1.Range("AQ58").Value2 = 1
$ws1.Range("AQ59").Value2 = 0
$ws1.Range("AQ60").Value2 = 6
$ws1.Range("AQ61").Value2 = 0
$ws1.Range("AQ62").Value2 = 1
$ws1.Range("AQ65").Value2 = 88
$ws1.Range("AQ66").Value2 = 282
$ws1.Range("AQ67").Value2 = 262
$ws1.Range("AQ68").Value2 = 5
$ws1.Range("AQ70").Value2 = 51
$ws1.Range("AQ71").Value2 = 31
$ws1.Range("AQ72").Value2 = 82
$ws1.Range("AQ73").Value2 = 22
$ws1.Range("AQ75").Value2 = 33
$ws1.Range("AQ76").Value2 = 74
$ws1.Range("AQ77").Value2 = 74
$ws1.Range("AQ78").Value2 = 2
$ws1.Range("AQ79").Value2 = 4

# --- Sheet "Total Cases by Ward": new date column R (4/17/2020) ---
$ws2 = $wb.Worksheets.Item("Total Cases by Ward")
# Normalize Q2 to the same date-number-format as the other header cells before adding R2
$ws2.Range("Q2").NumberFormat = $ws2.Range("P2").NumberFormat
$ws2.Range("R2").Value2 = 43937
$ws2.Range("R2").NumberFormat = $ws2.Range("P2").NumberFormat
$ws2.Range("R3").Value2 = 282
$ws2.Range("R4").Value2 = 207
$ws2.Range("R5").Value2 = 176
$ws2.Range("R6").Value2 = 423
$ws2.Range("R7").Value2 = 328
$ws2.Range("R8").Value2 = 337
$ws2.Range("R9").Value2 = 382
$ws2.Range("R10").Value2 = 316
$ws2.Range("R11").Value2 = 25

# --- Sheet "Total Cases by Race": new date column M (4/17/2020) ---
$ws3 = $wb.Worksheets.Item("Total Cases by Race")
$ws3.Range("M2").Value2 = 43937
$ws3.Range("M2").NumberFormat = $ws3.Range("L2").NumberFormat
$ws3.Range("M4").Value2 = 2476
$ws3.Range("M5").Value2 = 383
$ws3.Range("M6").Value2 = 461
$ws3.Range("M7").Value2 = 1171
$ws3.Range("M8").Value2 = 36
$ws3.Range("M9").Value2 = 7
$ws3.Range("M10").Value2 = 4
$ws3.Range("M11").Value2 = 386
$ws3.Range("M12").Value2 = 28
$ws3.Range("M14").Value2 = 544
$ws3.Range("M15").Value2 = 403
$ws3.Range("M16").Value2 = 1522
$ws3.Range("M17").Value2 = 7

# --- Sheet "Lives Lost by Race": new date column M (4/17/2020) ---
$ws4 = $wb.Worksheets.Item("Lives Lost by Race")
$ws4.Range("M1").Value2 = 43937
$ws4.Range("M1").NumberFormat = $ws4.Range("L1").NumberFormat
$ws4.Range("M3").Value2 = 86
$ws4.Range("M4").Value2 = 2
$ws4.Range("M5").Value2 = 66
$ws4.Range("M6").Value2 = 8
$ws4.Range("M7").Value2 = 10
$ws4.Range("M8").Value2 = 0

Write-Output "Applied 4/17/2020 data update to all sheets."
